# edit.ps1 - apply the commit's changes via PowerPoint COM interop
#
# The source diff does two things to the package:
#   1. ppt/slides/slide16.xml  - the table's <a:tableStyleId> changes from
#      {E2B9BD47-986C-4E41-8B6C-8E6843C5281A} to
#      {AF7D8CE3-6A9B-4676-9ADB-C7289818D8E1}.
#   2. ppt/theme/theme1.xml and ppt/theme/theme2.xml swap their entire
#      contents (name + color scheme) - theme1 becomes "Integral" (the
#      colours currently living in theme2) and theme2 becomes
#      "Office Theme" (the colours currently living in theme1). The font
#      scheme and format scheme are identical between the two themes, so
#      the only real content difference is the 12 colour-scheme entries
#      (and the name labels, which the exposed object model does not let
#      us rename).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style fix - find the (only) table on the deck and re-apply it
#    with the new style id. Table styles are not a plain settable
#    property; PowerPoint exposes Table.ApplyStyle(id) for this.
# ---------------------------------------------------------------------
$oldStyleId = "{E2B9BD47-986C-4E41-8B6C-8E6843C5281A}"
$newStyleId = "{AF7D8CE3-6A9B-4676-9ADB-C7289818D8E1}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme colour swap - the deck's live design theme (the one behind
#    SlideMaster/NotesMaster/Slide.ThemeColorScheme) is currently the
#    "Integral" palette; the commit swaps it for the "Office Theme"
#    palette that used to sit in the (otherwise COM-unreachable) second
#    theme part. Re-point every colour-scheme slot via
#    ColorScheme.Item(n).RGB - PowerPoint's RGB is packed as
#    0x00BBGGRR, so convert from the target hex (RRGGBB) accordingly.
# ---------------------------------------------------------------------
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# index -> target "Office Theme" colour (RRGGBB), in
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1..6, hlink, folHlink
$colorScheme.Item(1).RGB = 0          # dk1     000000
$colorScheme.Item(2).RGB = 16777215   # lt1     FFFFFF
$colorScheme.Item(3).RGB = 6968388    # dk2     44546A
$colorScheme.Item(4).RGB = 15132391   # lt2     E7E6E6
$colorScheme.Item(5).RGB = 13998939   # accent1 5B9BD5
$colorScheme.Item(6).RGB = 3243501    # accent2 ED7D31
$colorScheme.Item(7).RGB = 10855845   # accent3 A5A5A5
$colorScheme.Item(8).RGB = 49407      # accent4 FFC000
$colorScheme.Item(9).RGB = 12874308   # accent5 4472C4
$colorScheme.Item(10).RGB = 4697456   # accent6 70AD47
$colorScheme.Item(11).RGB = 12673797  # hlink   0563C1
$colorScheme.Item(12).RGB = 7491477   # folHlink 954F72
